# Update User Stories and Timeline
# Clear the now-removed "Waiting time of line" user-story rows and
# scroll/select the view to match the refreshed timeline.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F28").ClearContents()
$ws.Range("G30").ClearContents()
$ws.Range("B32:C33").ClearContents()

$ws.Range("C33").Select()
